$wb = $excel.ActiveWorkbook

# ---- Sheet: Hot Potato ----
$ws = $wb.Worksheets.Item("Hot Potato")
$ws.Range("A82:N82").Copy($ws.Range("A83:N83"))
$ws.Cells.Item(83, 1).Value = "ALLI"
$ws.Cells.Item(83, 2).Value = "R-T"
$ws.Cells.Item(83, 3).Value = "LILY"
$ws.Cells.Item(83, 4).Value = "DOUG"
$ws.Cells.Item(83, 5).Value = "MICO"
$ws.Cells.Item(83, 6).Value = "EMZ"
$ws.Cells.Item(83, 7).Value = "Equipo 2"
$ws.Cells.Item(83, 8).Value = "KCP|Fade"
$ws.Cells.Item(83, 9).Value = "KCP|Tyrant"
$ws.Cells.Item(83, 10).Value = "KCP|Zoulan"
$ws.Cells.Item(83, 11).Value = "SSG|bobby"
$ws.Cells.Item(83, 12).Value = "CODE|OG"
$ws.Cells.Item(83, 13).Value = "SSG|chino"
$ws.Cells.Item(83, 14).Value = "20250724T224249.000Z"
$gCell = $ws.Cells.Item(83, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true
$ws.Range("A83:N83").Copy($ws.Range("A84:N84"))
$ws.Cells.Item(84, 1).Value = "ALLI"
$ws.Cells.Item(84, 2).Value = "R-T"
$ws.Cells.Item(84, 3).Value = "LILY"
$ws.Cells.Item(84, 4).Value = "DOUG"
$ws.Cells.Item(84, 5).Value = "MICO"
$ws.Cells.Item(84, 6).Value = "EMZ"
$ws.Cells.Item(84, 7).Value = "Equipo 2"
$ws.Cells.Item(84, 8).Value = "KCP|Fade"
$ws.Cells.Item(84, 9).Value = "KCP|Tyrant"
$ws.Cells.Item(84, 10).Value = "KCP|Zoulan"
$ws.Cells.Item(84, 11).Value = "SSG|bobby"
$ws.Cells.Item(84, 12).Value = "CODE|OG"
$ws.Cells.Item(84, 13).Value = "SSG|chino"
$ws.Cells.Item(84, 14).Value = "20250724T224058.000Z"
$gCell = $ws.Cells.Item(84, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true

# ---- Sheet: Layer Cake ----
$ws = $wb.Worksheets.Item("Layer Cake")
$ws.Range("A76:N76").Copy($ws.Range("A77:N77"))
$ws.Cells.Item(77, 1).Value = "MR. P"
$ws.Cells.Item(77, 2).Value = "LOU"
$ws.Cells.Item(77, 3).Value = "KAZE"
$ws.Cells.Item(77, 4).Value = "DOUG"
$ws.Cells.Item(77, 5).Value = "GENE"
$ws.Cells.Item(77, 6).Value = "KENJI"
$ws.Cells.Item(77, 7).Value = "Equipo 1"
$ws.Cells.Item(77, 8).Value = "TE|Ezlivi"
$ws.Cells.Item(77, 9).Value = "TE|Rafikii"
$ws.Cells.Item(77, 10).Value = "TE|Belal"
$ws.Cells.Item(77, 11).Value = "TRB|Killer_17"
$ws.Cells.Item(77, 12).Value = "TRB|Zeus 解開"
$ws.Cells.Item(77, 13).Value = "TRB|R B M"
$ws.Cells.Item(77, 14).Value = "20250724T230443.000Z"
$gCell = $ws.Cells.Item(77, 7)
$gCell.Interior.Color = 16770508
$gCell.Font.Bold = $true
$ws.Range("A77:N77").Copy($ws.Range("A78:N78"))
$ws.Cells.Item(78, 1).Value = "MR. P"
$ws.Cells.Item(78, 2).Value = "LOU"
$ws.Cells.Item(78, 3).Value = "KAZE"
$ws.Cells.Item(78, 4).Value = "DOUG"
$ws.Cells.Item(78, 5).Value = "GENE"
$ws.Cells.Item(78, 6).Value = "KENJI"
$ws.Cells.Item(78, 7).Value = "Equipo 2"
$ws.Cells.Item(78, 8).Value = "TE|Ezlivi"
$ws.Cells.Item(78, 9).Value = "TE|Rafikii"
$ws.Cells.Item(78, 10).Value = "TE|Belal"
$ws.Cells.Item(78, 11).Value = "TRB|Killer_17"
$ws.Cells.Item(78, 12).Value = "TRB|Zeus 解開"
$ws.Cells.Item(78, 13).Value = "TRB|R B M"
$ws.Cells.Item(78, 14).Value = "20250724T230223.000Z"
$gCell = $ws.Cells.Item(78, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true
$ws.Range("A78:N78").Copy($ws.Range("A79:N79"))
$ws.Cells.Item(79, 1).Value = "CARL"
$ws.Cells.Item(79, 2).Value = "GENE"
$ws.Cells.Item(79, 3).Value = "GRAY"
$ws.Cells.Item(79, 4).Value = "GUS"
$ws.Cells.Item(79, 5).Value = "WILLOW"
$ws.Cells.Item(79, 6).Value = "LILY"
$ws.Cells.Item(79, 7).Value = "Equipo 2"
$ws.Cells.Item(79, 8).Value = "TE|Rafikii"
$ws.Cells.Item(79, 9).Value = "TE|Belal"
$ws.Cells.Item(79, 10).Value = "TE|Ezlivi"
$ws.Cells.Item(79, 11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(79, 12).Value = "TRB|R B M"
$ws.Cells.Item(79, 13).Value = "TRB|Lxffy"
$ws.Cells.Item(79, 14).Value = "20250724T225504.000Z"
$gCell = $ws.Cells.Item(79, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true
$ws.Range("A79:N79").Copy($ws.Range("A80:N80"))
$ws.Cells.Item(80, 1).Value = "CARL"
$ws.Cells.Item(80, 2).Value = "GENE"
$ws.Cells.Item(80, 3).Value = "GRAY"
$ws.Cells.Item(80, 4).Value = "GUS"
$ws.Cells.Item(80, 5).Value = "WILLOW"
$ws.Cells.Item(80, 6).Value = "LILY"
$ws.Cells.Item(80, 7).Value = "Equipo 2"
$ws.Cells.Item(80, 8).Value = "TE|Rafikii"
$ws.Cells.Item(80, 9).Value = "TE|Belal"
$ws.Cells.Item(80, 10).Value = "TE|Ezlivi"
$ws.Cells.Item(80, 11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(80, 12).Value = "TRB|R B M"
$ws.Cells.Item(80, 13).Value = "TRB|Lxffy"
$ws.Cells.Item(80, 14).Value = "20250724T225248.000Z"
$gCell = $ws.Cells.Item(80, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true
$ws.Range("A80:N80").Copy($ws.Range("A81:N81"))
$ws.Cells.Item(81, 1).Value = "MR. P"
$ws.Cells.Item(81, 2).Value = "LOU"
$ws.Cells.Item(81, 3).Value = "KAZE"
$ws.Cells.Item(81, 4).Value = "DOUG"
$ws.Cells.Item(81, 5).Value = "GENE"
$ws.Cells.Item(81, 6).Value = "KENJI"
$ws.Cells.Item(81, 7).Value = "Equipo 2"
$ws.Cells.Item(81, 8).Value = "TE|Ezlivi"
$ws.Cells.Item(81, 9).Value = "TE|Rafikii"
$ws.Cells.Item(81, 10).Value = "TE|Belal"
$ws.Cells.Item(81, 11).Value = "TRB|Killer_17"
$ws.Cells.Item(81, 12).Value = "TRB|Zeus 解開"
$ws.Cells.Item(81, 13).Value = "TRB|R B M"
$ws.Cells.Item(81, 14).Value = "20250724T230703.000Z"
$gCell = $ws.Cells.Item(81, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true
$ws.Range("A81:N81").Copy($ws.Range("A82:N82"))
$ws.Cells.Item(82, 1).Value = "HANK"
$ws.Cells.Item(82, 2).Value = "CHESTER"
$ws.Cells.Item(82, 3).Value = "MEEPLE"
$ws.Cells.Item(82, 4).Value = "CARL"
$ws.Cells.Item(82, 5).Value = "DOUG"
$ws.Cells.Item(82, 6).Value = "GRAY"
$ws.Cells.Item(82, 7).Value = "Equipo 2"
$ws.Cells.Item(82, 8).Value = "KCP|Fade"
$ws.Cells.Item(82, 9).Value = "KCP|Tyrant"
$ws.Cells.Item(82, 10).Value = "KCP|Zoulan"
$ws.Cells.Item(82, 11).Value = "SSG|chino"
$ws.Cells.Item(82, 12).Value = "CODE|OG"
$ws.Cells.Item(82, 13).Value = "SSG|bobby"
$ws.Cells.Item(82, 14).Value = "20250724T225419.000Z"
$gCell = $ws.Cells.Item(82, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true
$ws.Range("A82:N82").Copy($ws.Range("A83:N83"))
$ws.Cells.Item(83, 1).Value = "HANK"
$ws.Cells.Item(83, 2).Value = "CHESTER"
$ws.Cells.Item(83, 3).Value = "MEEPLE"
$ws.Cells.Item(83, 4).Value = "CARL"
$ws.Cells.Item(83, 5).Value = "DOUG"
$ws.Cells.Item(83, 6).Value = "GRAY"
$ws.Cells.Item(83, 7).Value = "Equipo 1"
$ws.Cells.Item(83, 8).Value = "KCP|Fade"
$ws.Cells.Item(83, 9).Value = "KCP|Tyrant"
$ws.Cells.Item(83, 10).Value = "KCP|Zoulan"
$ws.Cells.Item(83, 11).Value = "SSG|chino"
$ws.Cells.Item(83, 12).Value = "CODE|OG"
$ws.Cells.Item(83, 13).Value = "SSG|bobby"
$ws.Cells.Item(83, 14).Value = "20250724T225200.000Z"
$gCell = $ws.Cells.Item(83, 7)
$gCell.Interior.Color = 16770508
$gCell.Font.Bold = $true
$ws.Range("A83:N83").Copy($ws.Range("A84:N84"))
$ws.Cells.Item(84, 1).Value = "HANK"
$ws.Cells.Item(84, 2).Value = "CHESTER"
$ws.Cells.Item(84, 3).Value = "MEEPLE"
$ws.Cells.Item(84, 4).Value = "CARL"
$ws.Cells.Item(84, 5).Value = "DOUG"
$ws.Cells.Item(84, 6).Value = "GRAY"
$ws.Cells.Item(84, 7).Value = "Equipo 2"
$ws.Cells.Item(84, 8).Value = "KCP|Fade"
$ws.Cells.Item(84, 9).Value = "KCP|Tyrant"
$ws.Cells.Item(84, 10).Value = "KCP|Zoulan"
$ws.Cells.Item(84, 11).Value = "SSG|chino"
$ws.Cells.Item(84, 12).Value = "CODE|OG"
$ws.Cells.Item(84, 13).Value = "SSG|bobby"
$ws.Cells.Item(84, 14).Value = "20250724T224944.000Z"
$gCell = $ws.Cells.Item(84, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true

# ---- Sheet: Open Business ----
$ws = $wb.Worksheets.Item("Open Business")
$ws.Range("A78:N78").Copy($ws.Range("A79:N79"))
$ws.Cells.Item(79, 1).Value = "JESSIE"
$ws.Cells.Item(79, 2).Value = "GUS"
$ws.Cells.Item(79, 3).Value = "MOE"
$ws.Cells.Item(79, 4).Value = "KIT"
$ws.Cells.Item(79, 5).Value = "LOU"
$ws.Cells.Item(79, 6).Value = "CORDELIUS"
$ws.Cells.Item(79, 7).Value = "Equipo 1"
$ws.Cells.Item(79, 8).Value = "TE|Rafikii"
$ws.Cells.Item(79, 9).Value = "TE|Ezlivi"
$ws.Cells.Item(79, 10).Value = "TE|Belal"
$ws.Cells.Item(79, 11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(79, 12).Value = "TRB|R B M"
$ws.Cells.Item(79, 13).Value = "TRB|Lxffy"
$ws.Cells.Item(79, 14).Value = "20250724T224522.000Z"
$gCell = $ws.Cells.Item(79, 7)
$gCell.Interior.Color = 16770508
$gCell.Font.Bold = $true
$ws.Range("A79:N79").Copy($ws.Range("A80:N80"))
$ws.Cells.Item(80, 1).Value = "JESSIE"
$ws.Cells.Item(80, 2).Value = "GUS"
$ws.Cells.Item(80, 3).Value = "MOE"
$ws.Cells.Item(80, 4).Value = "KIT"
$ws.Cells.Item(80, 5).Value = "LOU"
$ws.Cells.Item(80, 6).Value = "CORDELIUS"
$ws.Cells.Item(80, 7).Value = "Equipo 1"
$ws.Cells.Item(80, 8).Value = "TE|Rafikii"
$ws.Cells.Item(80, 9).Value = "TE|Ezlivi"
$ws.Cells.Item(80, 10).Value = "TE|Belal"
$ws.Cells.Item(80, 11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(80, 12).Value = "TRB|R B M"
$ws.Cells.Item(80, 13).Value = "TRB|Lxffy"
$ws.Cells.Item(80, 14).Value = "20250724T224317.000Z"
$gCell = $ws.Cells.Item(80, 7)
$gCell.Interior.Color = 16770508
$gCell.Font.Bold = $true
$ws.Range("A80:N80").Copy($ws.Range("A81:N81"))
$ws.Cells.Item(81, 1).Value = "JESSIE"
$ws.Cells.Item(81, 2).Value = "GUS"
$ws.Cells.Item(81, 3).Value = "MOE"
$ws.Cells.Item(81, 4).Value = "KIT"
$ws.Cells.Item(81, 5).Value = "LOU"
$ws.Cells.Item(81, 6).Value = "CORDELIUS"
$ws.Cells.Item(81, 7).Value = "Equipo 2"
$ws.Cells.Item(81, 8).Value = "TE|Rafikii"
$ws.Cells.Item(81, 9).Value = "TE|Ezlivi"
$ws.Cells.Item(81, 10).Value = "TE|Belal"
$ws.Cells.Item(81, 11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(81, 12).Value = "TRB|R B M"
$ws.Cells.Item(81, 13).Value = "TRB|Lxffy"
$ws.Cells.Item(81, 14).Value = "20250724T224058.000Z"
$gCell = $ws.Cells.Item(81, 7)
$gCell.Interior.Color = 13421812
$gCell.Font.Bold = $true
